$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 221 (existing rows 221-233 shift down to 223-235)
$ws.Rows.Item(221).Resize(2).Insert(-4121)

# New row 221 data
$ws.Cells.Item(221, 1).Value = 5
$ws.Cells.Item(221, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(221, 3).Value = "Maule"
$ws.Cells.Item(221, 4).Value = 44516
$ws.Cells.Item(221, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(221, 5).Value = 7
$ws.Cells.Item(221, 6).Value = "Fruta"
$ws.Cells.Item(221, 7).Value = 100102
$ws.Cells.Item(221, 8).Value = "Cítricos"
$ws.Cells.Item(221, 9).Value = 100102004
$ws.Cells.Item(221, 10).Value = "Mandarina"
$ws.Cells.Item(221, 11).Value = "Murcott"
$ws.Cells.Item(221, 12).Value = "Especial"
$ws.Cells.Item(221, 13).Value = 190
$ws.Cells.Item(221, 14).Value = 8000
$ws.Cells.Item(221, 15).Value = 8000
$ws.Cells.Item(221, 16).Value = 8000
$ws.Cells.Item(221, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(221, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(221, 19).Value = 444
$ws.Cells.Item(221, 20).Value = 18

# New row 222 data
$ws.Cells.Item(222, 1).Value = 5
$ws.Cells.Item(222, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(222, 3).Value = "Maule"
$ws.Cells.Item(222, 4).Value = 44516
$ws.Cells.Item(222, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(222, 5).Value = 7
$ws.Cells.Item(222, 6).Value = "Fruta"
$ws.Cells.Item(222, 7).Value = 100102
$ws.Cells.Item(222, 8).Value = "Cítricos"
$ws.Cells.Item(222, 9).Value = 100102004
$ws.Cells.Item(222, 10).Value = "Mandarina"
$ws.Cells.Item(222, 11).Value = "Murcott"
$ws.Cells.Item(222, 12).Value = "Primera"
$ws.Cells.Item(222, 13).Value = 250
$ws.Cells.Item(222, 14).Value = 6000
$ws.Cells.Item(222, 15).Value = 6000
$ws.Cells.Item(222, 16).Value = 6000
$ws.Cells.Item(222, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(222, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(222, 19).Value = 333
$ws.Cells.Item(222, 20).Value = 18
